$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume-change (E) columns for rows 2-51.
# Values are written as text (matching the source data, which stores these as
# plain strings rather than numbers). Cells whose new value would otherwise be
# auto-parsed as a number by Excel are prefixed with a leading apostrophe to
# force text entry and preserve the exact digits/trailing zeros shown in the
# source (e.g. "0.6000" must stay "0.6000", not become the number 0.6).
$ws.Range("D2").Value = "26.606.57"
$ws.Range("E2").Value = "  -7.12%  "
$ws.Range("D3").Value = "1.693.69"
$ws.Range("E3").Value = "  -5.61%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'219.79"
$ws.Range("E5").Value = "  -4.96%  "
$ws.Range("E6").Value = "  -13.39%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.2653"
$ws.Range("E8").Value = "  -3.99%  "
$ws.Range("D9").Value = "'22.11"
$ws.Range("E9").Value = "  -4.41%  "
$ws.Range("D10").Value = "'0.06326"
$ws.Range("E10").Value = "  -6.16%  "
$ws.Range("D11").Value = "'0.07358"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").Value = "1.699.12"
$ws.Range("E12").Value = "  -5.37%  "
$ws.Range("D13").Value = "'4.519"
$ws.Range("E13").Value = "  -5.58%  "
$ws.Range("D14").Value = "'0.5788"
$ws.Range("E14").Value = "  -5.71%  "
$ws.Range("D15").Value = "1.924.68"
$ws.Range("E15").Value = "  -5.53%  "
$ws.Range("D16").Value = "'0.000008523"
$ws.Range("E16").Value = "  -5.74%  "
$ws.Range("D17").Value = "'65.39"
$ws.Range("E17").Value = "  -13.10%  "
$ws.Range("D18").Value = "26.629.85"
$ws.Range("E18").Value = "  -6.97%  "
$ws.Range("D19").Value = "'4.989"
$ws.Range("E19").Value = "  -8.58%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'10.96"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("E22").Value = "  -10.99%  "
$ws.Range("D23").Value = "'6.261"
$ws.Range("E23").Value = "  -8.06%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'144.59"
$ws.Range("E25").Value = "  -5.56%  "
$ws.Range("D26").Value = "'7.483"
$ws.Range("D27").Value = "'0.1166"
$ws.Range("E27").Value = "  -7.31%  "
$ws.Range("D28").Value = "'15.81"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("D30").Value = "'0.05743"
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("D31").Value = "'1.342"
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("D32").Value = "'3.524"
$ws.Range("E32").Value = "  -6.84%  "
$ws.Range("D33").Value = "'3.508"
$ws.Range("E33").Value = "  -7.95%  "
$ws.Range("D34").Value = "'1.640"
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "'0.6000"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("D37").Value = "'2.360"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").Value = "'2.680"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'0.01620"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").Value = "1.099.75"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "'0.8589"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").Value = "'5.832"
$ws.Range("E42").Value = "  -8.97%  "
$ws.Range("D43").Value = "'1.006"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'99.43"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "1.852.21"
$ws.Range("E45").Value = "  -4.91%  "
$ws.Range("E46").Value = "  +7.34%  "
$ws.Range("D47").Value = "'56.57"
$ws.Range("E47").Value = "  -5.49%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "'8.094"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "'0.4325"
$ws.Range("D51").Value = "'0.05235"
$ws.Range("E51").Value = "  -4.50%  "
